$p = $ppt.ActivePresentation
$m = $p.SlideMaster

# --- Shape 1 / id=2: Title Placeholder ----------------------------------
$shTitle = $m.Shapes.Item(1)
$shTitle.TextFrame.TextRange.Text = "Kliknutím lze upravit styl."

# --- Shape 2 / id=3: Text Placeholder -----------------------------------
$shBody = $m.Shapes.Item(2)
$shBody.TextFrame.TextRange.Paragraphs(1).Text = "Kliknutím lze upravit styly předlohy textu."
$shBody.TextFrame.TextRange.Paragraphs(2).Text = "Druhá úroveň"
$shBody.TextFrame.TextRange.Paragraphs(3).Text = "Třetí úroveň"
$shBody.TextFrame.TextRange.Paragraphs(4).Text = "Čtvrtá úroveň"
$shBody.TextFrame.TextRange.Paragraphs(5).Text = "Pátá úroveň"
